# Update XLSX file with new test cases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - ReCaptcha_1 / ReCaptcha component
$ws.Cells.Item(11, 1).Value = "ReCaptcha_1"
$ws.Cells.Item(11, 2).Value = "ReCaptcha component"
$ws.Cells.Item(11, 3).Value = "p1"
$ws.Cells.Item(11, 4).Value = "Verify that the CAPTCHA system prevents automated submissions or botting during the ticket purchasing phase."
$ws.Cells.Item(11, 5).Value = "CAPTCHA integrated on the purchase page"
$ws.Cells.Item(11, 6).Value = "1. Navigate to the ticket purchase page.`n2. Complete the CAPTCHA challenge incorrectly.`n3. Attempt to submit the purchase form."
$ws.Cells.Item(11, 7).Value = "CAPTCHA challenge must display an error message of failure when completed incorrectly. An error message`" (`"CAPTCHA validation failed`" is shown."
$ws.Cells.Item(11, 8).Value = "Error message `"CAPTCHA validation failed`" was displayed and the form was blocked."
$ws.Cells.Item(11, 9).Value = "Pass"
$ws.Cells.Item(11, 10).Value = "Carlyne"

# Row 12 - AdminPage_1 / Admin Page
$ws.Cells.Item(12, 1).Value = "AdminPage_1"
$ws.Cells.Item(12, 2).Value = "Admin Page"
$ws.Cells.Item(12, 3).Value = "p2"
$ws.Cells.Item(12, 4).Value = "Verify that the administrator can access all system data and perform actions on the page."
$ws.Cells.Item(12, 5).Value = "Admin credentials are available."
$ws.Cells.Item(12, 6).Value = "1. Log in using valid administrator credentials.`n2. Navigate to the admin page.`n3. Verify that all data (tickets, feedback, etc.) is displayed without filters."
$ws.Cells.Item(12, 7).Value = "Admin page should display system data with full access to required modifiers and data."
$ws.Cells.Item(12, 8).Value = "Admin login succeeded and all system data would be correctly displayed without filters that are shown with normal users."
$ws.Cells.Item(12, 9).Value = "Pass"
$ws.Cells.Item(12, 10).Value = "Carlyne"

# Row 13 - SessionTimeout_1 / Session Management
$ws.Cells.Item(13, 1).Value = "SessionTimeout_1"
$ws.Cells.Item(13, 2).Value = "Session Management"
$ws.Cells.Item(13, 3).Value = "p2"
$ws.Cells.Item(13, 4).Value = "Verify that a user is automatically logged out after 10 minutes of inactivity."
$ws.Cells.Item(13, 5).Value = "User is signed in and has an active session running."
$ws.Cells.Item(13, 6).Value = "1. Sign in and remain inactive for over 10 minutes.`n2. Attempt any interaction (ex: navigating to another page)."
$ws.Cells.Item(13, 7).Value = "The system should terminate the session and prompt the user to relog with a `"Session timed out`" notice."
$ws.Cells.Item(13, 8).Value = "The system terminated the session after 10 minutes of inactivity and redirected the user to the login page with a timeout message."
$ws.Cells.Item(13, 9).Value = "Pass"
$ws.Cells.Item(13, 10).Value = "Carlyne"

# Update the active selection to match the saved view state (G12)
$ws.Range("G12").Select()
